$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.065.45'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.75%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.834.29'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.28%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9995'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '239.85'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -2.00%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6713'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -3.20%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("B8").Value = 'Dogecoin'
$ws.Range("C8").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07439'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -3.32%  '
$ws.Range("B9").Value = 'Cardano'
$ws.Range("C9").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2969'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -2.89%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '22.94'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -3.30%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07658'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -1.47%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.831.13'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -1.44%  '
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -2.67%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6739'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -2.62%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '86.35'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -5.71%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.144'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -6.43%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '29.069.13'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -1.70%  '
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.55%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '227.23'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -5.37%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.47'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -2.35%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9995'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.05%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.315'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -3.73%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9997'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.02%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '160.26'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.36%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1429'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -4.62%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.681'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -2.70%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.74%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.506'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -1.75%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.236'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.34%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.119'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -1.50%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.197'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.26%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.05384'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +5.75%  '
$ws.Range("B33").Value = 'LidoDAOToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.857'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -1.94%  '
$ws.Range("B34").Value = 'ImmutableX'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7493'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -2.91%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -2.38%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.681'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.16%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.299.16'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -2.77%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01803'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -3.59%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.711'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.56%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9287'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -4.30%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.097'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +5.26%  '
$ws.Range("B42").Value = 'BabyDogeCoin'
$ws.Range("C42").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.00000000133'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +8.25%  '
$ws.Range("B43").Value = 'XinFinNetwork'
$ws.Range("C43").Value = 'https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.08391'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +30.07%  '
$ws.Range("B44").Value = 'Quant'
$ws.Range("C44").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '104.29'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -2.04%  '
$ws.Range("B45").Value = 'PaxDollar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.9988'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.11%  '
$ws.Range("B46").Value = 'RocketPoolETH'
$ws.Range("C46").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.973.75'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -1.34%  '
$ws.Range("B47").Value = 'Mantle'
$ws.Range("C47").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5175'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.78%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.447'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -3.26%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.756'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -1.10%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '63.54'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.09%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05928'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.09%  '
